# Apply "Penalty Reward System" edit to the PO data workbook.
#
# Sheet "Weekly Quantity" (first sheet): remove the data rows for
# order-weeks 2023-07-09 (45116.99999999999 / qty 690) and
# 2023-07-16 (45123.99999999999 / qty 690), and remove the row for
# order-week 2023-07-30 (45137.99999999999 / qty 890). The remaining
# row for order-week 2023-07-23 (45130.99999999999) has its
# "Requested quantity" changed from 950 to 420.
#
# Sheet "Monthly Trend" (second sheet): the "Requested quantity" for
# the month row dated 45138.99999999999 changes from 3220 to 420.

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

# --- Weekly Quantity sheet -------------------------------------------------

# Locate the rows to delete/edit by scanning column A for the known
# serial-date values, so the script is resilient to the exact row
# numbers in the sheet.
$lastRow = $wsWeekly.Cells.Item($wsWeekly.Rows.Count, 1).End(-4162).Row

$targetDates = @(45116.99999999999, 45123.99999999999, 45137.99999999999)
$rowsToDelete = @()
$editRow = $null

for ($r = $lastRow; $r -ge 2; $r--) {
    $val = $wsWeekly.Cells.Item($r, 1).Value2
    if ($null -eq $val) { continue }

    foreach ($d in $targetDates) {
        if ([math]::Abs($val - $d) -lt 0.0000001) {
            $rowsToDelete += $r
        }
    }

    if ([math]::Abs($val - 45130.99999999999) -lt 0.0000001) {
        $editRow = $r
    }
}

# Update the surviving row's quantity before removing the other rows.
if ($null -ne $editRow) {
    $wsWeekly.Cells.Item($editRow, 2).Value2 = 420
}

# Delete from the bottom up so row numbers of the other targets stay valid.
$rowsToDelete = $rowsToDelete | Sort-Object -Descending
foreach ($r in $rowsToDelete) {
    $wsWeekly.Rows.Item($r).Delete()
}

# --- Monthly Trend sheet ----------------------------------------------------

$lastRowMonthly = $wsMonthly.Cells.Item($wsMonthly.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRowMonthly; $r++) {
    $val = $wsMonthly.Cells.Item($r, 1).Value2
    if ($null -eq $val) { continue }
    if ([math]::Abs($val - 45138.99999999999) -lt 0.0000001) {
        $wsMonthly.Cells.Item($r, 2).Value2 = 420
        break
    }
}
